# Fruta / hortaliza, semanal
# Insert a new weekly record at row 46 (pushing the existing rows 46-60
# down to 47-61) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46:60 down to 47:61, creating a blank row 46.
$ws.Rows(46).Insert()

# Populate the new row 46 with the new weekly observation.
$ws.Cells.Item(46, 1).Value = 3
$ws.Cells.Item(46, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(46, 3).Value = "Coquimbo"
$ws.Cells.Item(46, 4).Value = 44466
$ws.Cells.Item(46, 5).Value = 5
$ws.Cells.Item(46, 6).Value = 100112026
$ws.Cells.Item(46, 7).Value = "Haba"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 38
$ws.Cells.Item(46, 11).Value = 13000
$ws.Cells.Item(46, 12).Value = 13000
$ws.Cells.Item(46, 13).Value = 13000
$ws.Cells.Item(46, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(46, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 16).Value = 520
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
